$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dialogs")
$ws.Activate()

$ws.Range("B2:B12").ClearContents()
$ws.Range("E9").Select()
